$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-27 Monday", "2024-05-28 Tuesday"),
    @("34×27=", "51×34="),
    @("85×54=", "42×59="),
    @("35×72=", "51×18="),
    @("12×11=", "69×36="),
    @("37×33=", "23×76="),
    @("27×17=", "86×67="),
    @("31×73=", "89×96="),
    @("31×30=", "90×58="),
    @("88×11=", "58×17="),
    @("89×35=", "72×30="),
    @("69×96=", "34×14="),
    @("90×42=", "15×92="),
    @("12×40=", "69×76="),
    @("95×32=", "62×65="),
    @("65×28=", "73×44="),
    @("14×25=", "18×20="),
    @("63×92=", "46×80="),
    @("57×27=", "92×46="),
    @("81×40=", "60×24="),
    @("97×51=", "27×66="),
    @("16×33=", "57×70="),
    @("47×59=", "26×44="),
    @("37×19=", "53×67="),
    @("74×91=", "32×75="),
    @("13×75=", "99×97=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
